$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Rewrite the first bullet title paragraph, splitting it into two
#    runs ("Pengembangan Kakas " + new title text) and move the
#    "_GoBack" bookmark onto it (it currently sits near the end of the
#    document, right after the big table).
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r = $p2.Range
[void]$r.MoveEnd(1, -1)   # exclude the paragraph mark
$r.Text = "Pengembangan Kakas Visualisasi Berbasis Web Untuk Graf dengan Modifikasi Matrix Framework"

$full = $p2.Range.Text
$splitIdx = $full.IndexOf("Visualisasi")
$splitPos = $p2.Range.Start + $splitIdx
$endPos = $p2.Range.End - 1
$rB = $d.Range($splitPos, $endPos)
# Toggling a formatting property forces the run to split in two
# without altering the visible character formatting.
$rB.Bold = 1
$rB.Bold = 0

# ------------------------------------------------------------------
# 2. Delete the next two bullet paragraphs entirely (the "Web-Based
#    Tool" variants) - they are removed outright by the edit.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$delRange = $d.Range($p3.Range.Start, $p4.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 3. Drop the old "_GoBack" bookmark (after the table, near the end of
#    the document) before re-adding it on the rewritten title above.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-fetch paragraph 2 and place a collapsed "_GoBack" bookmark right
# after its text (and before the paragraph mark). A collapsed range
# built directly on the paragraph-mark boundary confuses the engine,
# so insert a throwaway marker character at that spot, bookmark the
# marker's position, then delete the marker again; the bookmark stays
# put (collapsed) where the marker used to be.
$p2 = $d.Paragraphs.Item(2)
$insertPos = $p2.Range.End - 1
$markerRange = $d.Range($insertPos, $insertPos)
$markerRange.InsertAfter("@")

$p2 = $d.Paragraphs.Item(2)
$markerIdx = $p2.Range.Text.IndexOf("@")
$markerPos = $p2.Range.Start + $markerIdx
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($markerPos, $markerPos + 1).Delete()
